$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "57.756.06"
Set-TextValue "E2" "  +2.33%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.052.04"
Set-TextValue "E3" "  +2.27%  "

# Row 4 - TetherUSD
Set-TextValue "E4" "  -0.01%  "

# Row 5 - BNB
Set-TextValue "D5" "525.26"
Set-TextValue "E5" "  +5.83%  "

# Row 6 - Solana
Set-TextValue "D6" "142.41"
Set-TextValue "E6" "  +5.53%  "

# Row 7 - USDC
Set-TextValue "E7" "  -0.02%  "

# Row 8 - XRP
Set-TextValue "D8" "0.448"
Set-TextValue "E8" "  +5.03%  "

# Row 9 - Toncoin
Set-TextValue "D9" "7.64"
Set-TextValue "E9" "  +5.44%  "

# Row 10 - Dogecoin
Set-TextValue "E10" "  +7.74%  "

# Row 11 - Cardano
Set-TextValue "E11" "  +5.68%  "

# Row 12 - TRON
Set-TextValue "E12" "  +2.30%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue "D13" "3.576.44"
Set-TextValue "E13" "  +2.29%  "

# Row 14 - Avalanche
Set-TextValue "D14" "27.04"
Set-TextValue "E14" "  +8.14%  "

# Row 15 - ShibaInu
Set-TextValue "E15" "  +16.90%  "

# Row 16 - Polkadot
Set-TextValue "D16" "6.30"
Set-TextValue "E16" "  +8.06%  "

# Row 17 - WrappedBTC
Set-TextValue "D17" "57.713.84"
Set-TextValue "E17" "  +2.43%  "

# Row 18 - WrappedEther
Set-TextValue "D18" "3.061.22"
Set-TextValue "E18" "  +2.34%  "

# Row 19 - Chainlink
Set-TextValue "D19" "13.09"
Set-TextValue "E19" "  +6.01%  "

# Row 20 - Uniswap
Set-TextValue "D20" "8.18"
Set-TextValue "E20" "  +5.64%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "337.94"
Set-TextValue "E21" "  +3.84%  "

# Row 22 - Dai
Set-TextValue "E22" "  +0.07%  "

# Row 23 - Polygon
Set-TextValue "D23" "0.502"
Set-TextValue "E23" "  +7.40%  "

# Row 24 - Litecoin
Set-TextValue "D24" "64.93"
Set-TextValue "E24" "  +5.91%  "

# Row 25 - Kaspa
Set-TextValue "E25" "  +6.93%  "

# Row 26 - PEPE
Set-TextValue "E26" "  +9.10%  "

# Row 27 - Binance-PegBSC-USD
Set-TextValue "D27" "0.999"
Set-TextValue "E27" "  +0.18%  "

# Row 28 - RenderToken
Set-TextValue "D28" "6.92"
Set-TextValue "E28" "  +6.23%  "

# Row 29 - InternetComputer(DFINITY)
Set-TextValue "D29" "7.38"
Set-TextValue "E29" "  +10.80%  "

# Row 30 - PancakeSwap
Set-TextValue "D30" "1.86"
Set-TextValue "E30" "  +6.54%  "

# Row 31 - Fetch.AI
Set-TextValue "E31" "  +5.32%  "

# Row 32 - EthereumClassic
Set-TextValue "D32" "21.14"
Set-TextValue "E32" "  +4.62%  "

# Row 33 - Monero
Set-TextValue "D33" "156.69"
Set-TextValue "E33" "  +1.06%  "

# Row 34 - NEARProtocol
Set-TextValue "E34" "  +6.28%  "

# Row 35 - Aptos
Set-TextValue "D35" "6.01"
Set-TextValue "E35" "  +7.25%  "

# Row 36 - ImmutableX
Set-TextValue "E36" "  +3.36%  "

# Row 37 - EnergySwap
Set-TextValue "D37" "26.16"
Set-TextValue "E37" "  +12.91%  "

# Row 38 - Hedera
Set-TextValue "D38" "0.0705"
Set-TextValue "E38" "  +2.78%  "

# Row 39 - RenzoRestakedETH
Set-TextValue "D39" "3.088.14"
Set-TextValue "E39" "  +2.28%  "

# Row 40 - OKB
Set-TextValue "D40" "37.73"
Set-TextValue "E40" "  +3.23%  "

# Row 41 - Filecoin
Set-TextValue "D41" "3.90"
Set-TextValue "E41" "  +9.16%  "

# Row 42 - FirstDigitalUSD
Set-TextValue "E42" "  +0.02%  "

# Row 43 - was Stacks, now Mantle
Set-TextValue "B43" "Mantle"
Set-TextValue "C43" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D43" "0.664"
Set-TextValue "E43" "  +3.72%  "

# Row 44 - was Mantle, now Stacks
Set-TextValue "B44" "Stacks"
Set-TextValue "C44" "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D44" "1.47"
Set-TextValue "E44" "  +5.21%  "

# Row 45 - Maker
Set-TextValue "D45" "2.327.88"
Set-TextValue "E45" "  +4.25%  "

# Row 46 - ONDO
Set-TextValue "E46" "  +3.93%  "

# Row 47 - dogwifhat
Set-TextValue "E47" "  +3.80%  "

# Row 48 - VeChain
Set-TextValue "E48" "  +4.16%  "

# Row 49 - Cosmos
Set-TextValue "E49" "  +4.49%  "

# Row 50 - InjectiveProtocol
Set-TextValue "E50" "  +5.98%  "

# Row 51 - Stellar
Set-TextValue "E51" "  +6.27%  "
